# quarterly.xlsx update — "update database and change read_price algorithm"
#
# The quarterly report rolls its 10-quarter window forward by one quarter:
#   - drops "فصل دوم منتهی به 1399/06" (Q2 1399)
#   - appends "فصل چهارم منتهی به 1401/12" (Q4 1401)
# All figures in the two data tables (columns E:N) shift one column to the
# left (E<-F, F<-G, ... M<-N) and the new quarter's figures land in column N.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("E", "F", "G", "H", "I", "J", "K", "L", "M", "N")

# ---------------------------------------------------------------------
# 1. Header labels (row 8 and row 24) — rolling quarter window
# ---------------------------------------------------------------------
$quarters = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "8").Value = $quarters[$i]
    $ws.Range($cols[$i] + "24").Value = $quarters[$i]
}

# ---------------------------------------------------------------------
# 2. First table (رows 10-20) — هزینه های عمومی و اداری
# ---------------------------------------------------------------------
$rowValues = @{
    10 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 357543)
    12 = @(0, 0, 0, 0, 0, 11006, 2490, 2510, -2719, 760)
    13 = @(0, 0, 0, 0, 0, 0, 0, 0, 2241, 747)
    14 = @(367, 2495, 2500, 290, 5387, 4198, 2812, 3357, -6169, 0)
    15 = @(835, -305, 250, 518, 480, 1059, 842, 1155, -1717, 93)
    16 = @(1537, -155, 679, 805, 742, 382, 828, 1063, 6059, 2650)
    17 = @(11184, 20965, 18826, 15871, 17349, 23093, 27493, 20871, 34090, 27485)
    19 = @(6854, 13759, 18912, 1115, -2443, 22301, 12200, 21783, 10903, 14962)
    20 = @(20777, 36759, 41167, 18599, 21515, 62039, 46665, 50739, 42688, 404240)
}

foreach ($r in $rowValues.Keys) {
    $vals = $rowValues[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}

# ---------------------------------------------------------------------
# 3. Second table (rows 26-27) — تعداد پرسنل
# ---------------------------------------------------------------------
$personnelValues = @{
    26 = @(46, 35, 47, 35, 35, 23, 28, 29, 27, 27)
    27 = @(182, 197, 185, 197, 197, 204, 200, 200, 209, 212)
}

foreach ($r in $personnelValues.Keys) {
    $vals = $personnelValues[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}
